# Add two new transaction rows (43 and 44) to the "Наличные" (cash) sheet,
# mirroring the existing ledger rows for the same user/file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(43, 1).Value = 7676096317
$ws.Cells.Item(43, 2).Value = "M-Банкинг чек-4294968802.pdf"
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = 0.2
$ws.Cells.Item(43, 5).Value = "2025-06-28 20:15:42"

$ws.Cells.Item(44, 1).Value = 7676096317
$ws.Cells.Item(44, 2).Value = "M-Банкинг чек-4294968802.pdf"
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = 0.2
$ws.Cells.Item(44, 5).Value = "2025-06-28 20:16:28"
